$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Abstract section: split the final paragraph so the trailing (hidden)
#    _GoBack bookmark ends up alone in its own paragraph, and insert an
#    empty paragraph plus a new paragraph of text in between.
# ---------------------------------------------------------------------------

$lastPara = $d.Paragraphs.Last
$searchRange = $lastPara.Range.Duplicate
$found = $searchRange.Find.Execute( `
    "Aspektes soll jeweils in einer Iteration komplett analysiert, designt, implementiert und getestet werden.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$splitPoint = $searchRange.End
$splitRange = $d.Range($splitPoint, $splitPoint)

$cr = [string][char]13
$splitRange.InsertBefore($cr + $cr + $cr)

# Paragraph that will hold the new blank line.
$blankPara = $d.Paragraphs.Item($lastPara.Index + 1)
$blankFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:before="120"/><w:jc w:val="both"/><w:rPr><w:sz w:val="20"/></w:rPr></w:pPr></w:p>'
$blankPara.Range.InsertXML($blankFrag)

# Paragraph that will hold the new "Eclipse RCP..." text.
$textPara = $d.Paragraphs.Item($lastPara.Index + 2)
$textFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:before="120"/><w:jc w:val="both"/><w:rPr><w:sz w:val="20"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="20"/></w:rPr><w:t>Eclipse</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve"> RCP ist ein Standardframework f' + [string][char]0xFC + 'r Gesch' + [string][char]0xE4 + 'ftsanwendungen. Mit der neusten Generation E4 wurde </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="20"/></w:rPr><w:t>Eclipse</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve"> RCP vollst' + [string][char]0xE4 + 'ndig modernisiert. Anhand einer wichtigen RCP Applikation der SBB wird eine Migration auf </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="20"/></w:rPr><w:t>Eclipse</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve"> E4 exemplarisch durchgef' + [string][char]0xFC + 'hrt und die dabei ber' + [string][char]0xFC + 'cksichtigten Aspekte dargestellt.</w:t></w:r></w:p>'
$textPara.Range.InsertXML($textFrag)

# ---------------------------------------------------------------------------
# 2) Header logo: Word regenerated the drawing's anchorId/editId when the
#    image was touched again; reproduce that by re-stamping the drawing.
# ---------------------------------------------------------------------------

$hdr = $d.Sections.Item(1).Headers.Item(1)
$hdrRange = $hdr.Range
$headerFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" w:rsidR="00AB072F" w:rsidRDefault="0061617C" w:rsidP="00A14B2F"><w:pPr><w:jc w:val="center"/><w:rPr><w:sz w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/><w:sz w:val="20"/><w:lang w:val="de-CH"/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="22F1C139" wp14:editId="2F8E69B8"><wp:extent cx="2628900" cy="619125"/><wp:effectExtent l="19050" t="0" r="0" b="0"/><wp:docPr id="1" name="Bild 2"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="0" name="Bild 2"/><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId1"/><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="2628900" cy="619125"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln w="9525"><a:noFill/><a:miter lim="800000"/><a:headEnd/><a:tailEnd/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>'
$hdrRange.InsertXML($headerFrag)

Write-Output "done"
